$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 412 - this shifts the existing rows 412:500
# down to 413:501 (matching the diff, which shows every row from 412
# onward taking on the values previously held by the row above it, and
# a brand-new row appearing at the end as row 501).
$ws.Rows.Item(412).Insert()

# Populate the newly inserted row 412 with the new weekly price record.
$ws.Cells.Item(412, 1).Value2 = 9
$ws.Cells.Item(412, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(412, 3).Value2 = "Metropolitana"
$ws.Cells.Item(412, 4).Value2 = 45211
$ws.Cells.Item(412, 5).Value2 = 13
$ws.Cells.Item(412, 6).Value2 = 300000001
$ws.Cells.Item(412, 7).Value2 = "Rabanito"
$ws.Cells.Item(412, 8).Value2 = "Sin especificar"
$ws.Cells.Item(412, 9).Value2 = "Primera"
$ws.Cells.Item(412, 10).Value2 = 7000
$ws.Cells.Item(412, 11).Value2 = 3000
$ws.Cells.Item(412, 12).Value2 = 3000
$ws.Cells.Item(412, 13).Value2 = 3000
$ws.Cells.Item(412, 14).Value2 = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(412, 15).Value2 = "Provincia de Chacabuco"
$ws.Cells.Item(412, 16).Value2 = 30
$ws.Cells.Item(412, 17).Value2 = 100
$ws.Cells.Item(412, 18).Value2 = "Hortaliza"
